$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'60.501.91"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.44%  "

$ws.Range("D3").Value = "'2.605.07"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.17%  "

$ws.Range("E4").Value = "  +0.09%  "

$ws.Range("D5").Value = "'514.75"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.05%  "

$ws.Range("D6").Value = "'153.80"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.46%  "

$ws.Range("E7").Value = "  +0.20%  "

$ws.Range("D8").Value = "'0.600"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +4.52%  "

$ws.Range("D9").Value = "'2.616.54"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.42%  "

$ws.Range("E10").Value = "  +3.96%  "

$ws.Range("E11").Value = "  +0.38%  "

$ws.Range("D12").Value = "'0.346"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.75%  "

$ws.Range("E13").Value = "  +1.87%  "

$ws.Range("D14").Value = "'3.061.81"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.77%  "

$ws.Range("D15").Value = "'60.538.25"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.45%  "

$ws.Range("D16").Value = "'21.65"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.35%  "

$ws.Range("E17").Value = "  +1.34%  "

$ws.Range("D18").Value = "'2.607.82"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.12%  "

$ws.Range("E19").Value = "  -0.40%  "

$ws.Range("D20").Value = "'357.64"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +5.19%  "

$ws.Range("D21").Value = "'10.60"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.37%  "

$ws.Range("D22").Value = "'6.20"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.15%  "

$ws.Range("E23").Value = "  +0.24%  "

$ws.Range("D24").Value = "'61.06"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.72%  "

$ws.Range("E25").Value = "  +1.29%  "

$ws.Range("D26").Value = "'2.726.82"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.81%  "

$ws.Range("E27").Value = "  +0.70%  "

$ws.Range("D28").Value = "'0.994"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.08%  "

$ws.Range("E29").Value = "  -1.35%  "

$ws.Range("D30").Value = "'7.33"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.72%  "

$ws.Range("E31").Value = "  +0.19%  "

$ws.Range("E32").Value = "  +0.90%  "

$ws.Range("E33").Value = "  +1.78%  "

$ws.Range("B34").Value = "Aptos"
$ws.Range("C34").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D34").Value = "'5.90"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.43%  "

$ws.Range("B35").Value = "Monero"
$ws.Range("C35").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D35").Value = "'150.45"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.56%  "

$ws.Range("E36").Value = "  +1.36%  "

$ws.Range("E37").Value = "  -0.67%  "

$ws.Range("D38").Value = "'0.891"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.96%  "

$ws.Range("E39").Value = "  +0.76%  "

$ws.Range("D40").Value = "'0.845"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.51%  "

$ws.Range("D41").Value = "'36.21"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.04%  "

$ws.Range("D42").Value = "'3.74"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.66%  "

$ws.Range("D43").Value = "'289.87"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.57%  "

$ws.Range("E44").Value = "  +2.07%  "

$ws.Range("D45").Value = "'0.620"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.05%  "

$ws.Range("D46").Value = "'0.996"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.08%  "

$ws.Range("D47").Value = "'0.0556"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.20%  "

$ws.Range("E48").Value = "  -0.57%  "

$ws.Range("D49").Value = "'4.96"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.02%  "

$ws.Range("E50").Value = "  +1.43%  "

$ws.Range("E51").Value = "  +0.47%  "

